$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N (2021 data), rows 2-15 -----------------------------------

# Row 2: blank "border" cell -> copy format from M2 (style s=8)
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

# Row 3: header year 2021 -> copy format from M3 (style s=6)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2021

# Row 4: copy format from M4 (style s=15, bold percent) - set below after the
# D4:L4 re-style, so M4 already carries the right style at that point.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 95.134712433469176

# Row 5 & 6: copy format from D5/D6 (style s=10) - matches M5/M6 already.
$ws.Range("D5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 99.705541665880986

$ws.Range("D6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 92.425193326577897

# Rows 7-14: copy format from D<row> (style s=10) - note this differs from
# M<row>, which is right-aligned (s=11); the new column keeps the plain
# vertical-centered numeric style instead.
$ws.Range("D7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 88.209991167538519

$ws.Range("D8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 92.225038985690773

$ws.Range("D9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 96.801032063987265

$ws.Range("D10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 97.660491031729507

$ws.Range("D11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 90.23262877800066

$ws.Range("D12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 99.653994395099105

$ws.Range("D13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 100

$ws.Range("D14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 100

# Row 15: copy format from M15 (style s=13)
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = 100

# --- Row 4 (D4:L4) switches from the plain numeric style to the bold one ---
# already used by M4 (xf 10 -> xf 15, i.e. the same font/format but bold).
$ws.Range("D4:L4").Font.Bold = $true

# --- Selection moves to the newly entered cell ------------------------------
$ws.Range("N2").Select()

# --- Printer resolution bump (paper stays the same, DPI 0 -> 300) ----------
$ws.PageSetup.PrintQuality = 300
